$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 5000931.5
$ws.Range("I46").Value = 497
$ws.Range("J46").Value = 5435751.5
$ws.Range("K46").Value = 1491
$ws.Range("L46").Value = 16307254.5
$ws.Range("M46").Value = -1372
$ws.Range("N46").Value = -16307492.5
$ws.Range("H60").Value = 5000931.5
$ws.Range("I60").Value = 497
$ws.Range("J60").Value = 5435751.5
$ws.Range("K60").Value = 1491
$ws.Range("L60").Value = 16307254.5
$ws.Range("M60").Value = -1007
$ws.Range("N60").Value = -16308222.5
$ws.Range("H64").Value = 29902.297
$ws.Range("I64").Value = 69339.92999999999
$ws.Range("K64").Value = 69339.92999999999
$ws.Range("M64").Value = -69091.92999999999
$ws.Range("H67").Value = 29902.297
$ws.Range("I67").Value = 69339.92999999999
$ws.Range("K67").Value = 69339.92999999999
$ws.Range("M67").Value = -68481.92999999999
$ws.Range("H106").Value = 2000.5
$ws.Range("J106").Value = 2801.6667
$ws.Range("L106").Value = 2801.6667
$ws.Range("N106").Value = -4063.6667
$ws.Range("H113").Value = 3099.6667
$ws.Range("J113").Value = 1999
$ws.Range("L113").Value = 1999
$ws.Range("N113").Value = -8507
$ws.Range("H129").Value = 2336.242
$ws.Range("I129").Value = 7614.0713
$ws.Range("J129").Value = 796.875
$ws.Range("K129").Value = 22842.2139
$ws.Range("L129").Value = 2390.625
$ws.Range("M129").Value = -17842.2139
$ws.Range("N129").Value = -12390.625

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 14116.667
$ws.Range("I37").Value = 9800
$ws.Range("J37").Value = 14980
$ws.Range("K37").Value = 9800
$ws.Range("L37").Value = 14980
$ws.Range("M37").Value = -9527
$ws.Range("N37").Value = -15526
$ws.Range("H55").Value = 13025
$ws.Range("J55").Value = 13885.714
$ws.Range("L55").Value = 13885.714
$ws.Range("N55").Value = -14515.714
$ws.Range("H62").Value = 35249
$ws.Range("J62").Value = 35249
$ws.Range("L62").Value = 35249
$ws.Range("N62").Value = -36497
$ws.Range("H65").Value = 35249
$ws.Range("J65").Value = 35249
$ws.Range("L65").Value = 105747
$ws.Range("N65").Value = -111987
$ws.Range("H80").Value = 26528.25
$ws.Range("J80").Value = 27460.857
$ws.Range("L80").Value = 27460.857
$ws.Range("N80").Value = -29456.857
$ws.Range("H83").Value = 26528.25
$ws.Range("J83").Value = 27460.857
$ws.Range("L83").Value = 82382.571
$ws.Range("N83").Value = -92366.571
$ws.Range("H92").Value = 19000
$ws.Range("J92").Value = 19000
$ws.Range("L92").Value = 19000
$ws.Range("N92").Value = -23992
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
$ws.Range("H110").Value = 58888950
$ws.Range("I110").Value = 71507864
$ws.Range("J110").Value = 684.3333
$ws.Range("K110").Value = 71507864
$ws.Range("L110").Value = 684.3333
$ws.Range("M110").Value = -71505819
$ws.Range("N110").Value = -4774.3333

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 48560.332
$ws.Range("J62").Value = 48560.332
$ws.Range("L62").Value = 48560.332
$ws.Range("N62").Value = -49932.332
$ws.Range("H65").Value = 48560.332
$ws.Range("J65").Value = 48560.332
$ws.Range("L65").Value = 145680.996
$ws.Range("N65").Value = -152544.996
$ws.Range("H105").Value = 112996.555
$ws.Range("I105").Value = 126871.125
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 126871.125
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -125124.125
$ws.Range("N105").Value = -5494
$ws.Range("H112").Value = 36500
$ws.Range("J112").Value = 36500
$ws.Range("L112").Value = 36500
$ws.Range("N112").Value = -39454
$ws.Range("H134").Value = 2784
$ws.Range("I134").Value = 2285.1
$ws.Range("J134").Value = 4922.143
$ws.Range("K134").Value = 6855.299999999999
$ws.Range("L134").Value = 14766.429
$ws.Range("M134").Value = -4320.299999999999
$ws.Range("N134").Value = -19836.429

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3243.5789
$ws.Range("J31").Value = 4471.4
$ws.Range("L31").Value = 4471.4
$ws.Range("N31").Value = -5061.4
$ws.Range("H34").Value = 3243.5789
$ws.Range("J34").Value = 4471.4
$ws.Range("L34").Value = 4471.4
$ws.Range("N34").Value = -4875.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1248.8485
$ws.Range("J5").Value = 1155.174
$ws.Range("L5").Value = 3465.522
$ws.Range("N5").Value = -3689.522
$ws.Range("H87").Value = 8093.3335
$ws.Range("I87").Value = 7480
$ws.Range("K87").Value = 22440
$ws.Range("M87").Value = -21192
$ws.Range("H90").Value = 8093.3335
$ws.Range("I90").Value = 7480
$ws.Range("K90").Value = 67320
$ws.Range("M90").Value = -61080
$ws.Range("H100").Value = 3374.875
$ws.Range("J100").Value = 3374.875
$ws.Range("L100").Value = 10124.625
$ws.Range("N100").Value = -11746.625
$ws.Range("H113").Value = 832
$ws.Range("I113").Value = 1229.3077
$ws.Range("J113").Value = 545.05554
$ws.Range("K113").Value = 3687.9231
$ws.Range("L113").Value = 1635.16662
$ws.Range("M113").Value = -1517.9231
$ws.Range("N113").Value = -5975.16662
$ws.Range("H135").Value = 1248.8485
$ws.Range("J135").Value = 1155.174
$ws.Range("L135").Value = 10396.566
$ws.Range("N135").Value = -15466.566

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 844003.2
$ws.Range("I46").Value = 487.25
$ws.Range("J46").Value = 1265761.1
$ws.Range("K46").Value = 487.25
$ws.Range("L46").Value = 1265761.1
$ws.Range("M46").Value = -299.25
$ws.Range("N46").Value = -1266137.1
$ws.Range("H93").Value = 1242.2954
$ws.Range("I93").Value = 1145.8334
$ws.Range("J93").Value = 1449
$ws.Range("K93").Value = 1145.8334
$ws.Range("L93").Value = 1449
$ws.Range("M93").Value = 102.1666
$ws.Range("N93").Value = -3945
$ws.Range("H98").Value = 29490
$ws.Range("J98").Value = 29490
$ws.Range("L98").Value = 29490
$ws.Range("N98").Value = -35480
$ws.Range("H106").Value = 22083.334
$ws.Range("J106").Value = 22083.334
$ws.Range("L106").Value = 22083.334
$ws.Range("N106").Value = -24607.334
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
